$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mapping")

$ws.Range("A2").Value = -77.333
$ws.Range("B2").Value = -77.0652

$ws.Range("A3").Value = 37.936
$ws.Range("B3").Value = 38.1453

$ws.Range("A4").Value = -74.9823
$ws.Range("B4").Value = -75.2577

$ws.Range("A5").Value = 39.7394
$ws.Range("B5").Value = 39.532
